$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Pre-Noon" (rows 8-9) and "Post-Noon" (rows 10-11) blocks of this
# MultiIndex-style table need to swap places (the two row-groups trade
# positions), while the column headers (rows 5-6) stay put.
#
# Capture the current ("before") values for columns I:N across rows 8-11
# first (using Value2, which reads the real cell contents), then write
# them back swapped, so the read pass never sees a partially-updated
# state.

$colsTop = @("I8","J8","K8","L8","M8","N8")
$colsTopLower = @("I9","J9","K9","L9","M9","N9")
$colsBottom = @("I10","J10","K10","L10","M10","N10")
$colsBottomLower = @("I11","J11","K11","L11","M11","N11")

$topVals = @{}
foreach ($addr in $colsTop) { $topVals[$addr] = $ws.Range($addr).Value2 }
$topLowerVals = @{}
foreach ($addr in $colsTopLower) { $topLowerVals[$addr] = $ws.Range($addr).Value2 }
$bottomVals = @{}
foreach ($addr in $colsBottom) { $bottomVals[$addr] = $ws.Range($addr).Value2 }
$bottomLowerVals = @{}
foreach ($addr in $colsBottomLower) { $bottomLowerVals[$addr] = $ws.Range($addr).Value2 }

# Row 8 <- old Row 10 ; Row 10 <- old Row 8
for ($i = 0; $i -lt $colsTop.Count; $i++) {
    $ws.Range($colsTop[$i]).Value = $bottomVals[$colsBottom[$i]]
    $ws.Range($colsBottom[$i]).Value = $topVals[$colsTop[$i]]
}

# Row 9 <- old Row 11 ; Row 11 <- old Row 9
for ($i = 0; $i -lt $colsTopLower.Count; $i++) {
    $ws.Range($colsTopLower[$i]).Value = $bottomLowerVals[$colsBottomLower[$i]]
    $ws.Range($colsBottomLower[$i]).Value = $topLowerVals[$colsTopLower[$i]]
}
